$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.077.69"
$ws.Range("E2").Value = "  -2.91%  "
$ws.Range("D3").Value = "1.714.14"
$ws.Range("E3").Value = "  -3.21%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'308.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.95%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "'0.4646"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.52%  "
$ws.Range("D8").Value = "'0.3419"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.00%  "
$ws.Range("D9").Value = "'41.96"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("D10").Value = "'0.07247"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.90%  "
$ws.Range("D11").Value = "'1.042"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.33%  "
$ws.Range("D12").Value = "'1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "'19.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.51%  "
$ws.Range("D14").Value = "'5.845"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.18%  "
$ws.Range("D15").Value = "1.712.09"
$ws.Range("E15").Value = "  -3.24%  "
$ws.Range("D16").Value = "'6.860"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.27%  "
$ws.Range("D17").Value = "'88.62"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.08%  "
$ws.Range("D18").Value = "'0.00001036"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.24%  "
$ws.Range("D19").Value = "'0.06344"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.45%  "
$ws.Range("D21").Value = "'16.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.16%  "
$ws.Range("D22").Value = "'5.640"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.57%  "
$ws.Range("D23").Value = "27.107.05"
$ws.Range("E23").Value = "  -2.98%  "
$ws.Range("D24").Value = "'10.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.01%  "
$ws.Range("D25").Value = "'2.135"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.26%  "
$ws.Range("D26").Value = "'156.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.09%  "
$ws.Range("D27").Value = "'19.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.57%  "
$ws.Range("D28").Value = "1.912.74"
$ws.Range("E28").Value = "  -3.04%  "
$ws.Range("D29").Value = "'2.111"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.72%  "
$ws.Range("D30").Value = "'119.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.55%  "
$ws.Range("D31").Value = "'1.021"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.99%  "
$ws.Range("D32").Value = "'0.09135"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("D33").Value = "'3.601"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("D34").Value = "'5.319"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.10%  "
$ws.Range("D35").Value = "'0.02190"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.70%  "
$ws.Range("D36").Value = "'0.05813"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.12%  "
$ws.Range("D37").Value = "'11.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.32%  "
$ws.Range("D38").Value = "'0.1993"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.08%  "
$ws.Range("D39").Value = "'4.722"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.98%  "
$ws.Range("D40").Value = "'1.390"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("D41").Value = "'0.5902"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.62%  "
$ws.Range("D42").Value = "'1.124"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.94%  "
$ws.Range("D43").Value = "'7.446"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.36%  "
$ws.Range("D44").Value = "'12.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.17%  "
$ws.Range("D45").Value = "'3.558"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.94%  "
$ws.Range("D46").Value = "'0.5638"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.02%  "
$ws.Range("D47").Value = "'118.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.30%  "
$ws.Range("D48").Value = "'1.841"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.86%  "
$ws.Range("D49").Value = "'0.06648"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.76%  "
$ws.Range("D50").Value = "'1.080"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.05%  "
